$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.177.15"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.629.27"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.608.73"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "27.155.88"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").Value = "1.313.17"
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D43").Value = "1.767.64"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  -4.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "0.0₆0106"
$ws.Range("E47").Value = "  +6.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.805"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +20.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.91%  "
